$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.467281579971313
$ws.Range("B1").Value = 3.221961975097656
$ws.Range("C1").Value = 2.777711153030396
$ws.Range("D1").Value = 2.242633581161499
$ws.Range("E1").Value = 1.456994295120239
